$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price text would otherwise be auto-coerced to a number by Excel
# (single decimal point, no thousands separators) need NumberFormat forced to
# Text before the assignment, then restored so no stray style sticks around.

$ws.Range("D2").Value = "42.147.15"
$ws.Range("E2").Value = "  -1.99%  "

$ws.Range("D3").Value = "2.269.39"
$ws.Range("E3").Value = "  -3.00%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "297.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.17%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.76"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.76%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("E8").Value = "  -3.76%  "

$ws.Range("E9").Value = "  -3.91%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.32"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.46%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0788"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.25%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "48.16"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -8.25%  "

$ws.Range("E13").Value = "  -0.28%  "

$ws.Range("E14").Value = "  -3.24%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.60"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.57%  "

$ws.Range("D16").Value = "2.622.27"
$ws.Range("E16").Value = "  -3.02%  "

$ws.Range("D17").Value = "2.265.50"
$ws.Range("E17").Value = "  -3.51%  "

$ws.Range("E18").Value = "  -6.17%  "

$ws.Range("D19").Value = "42.084.51"
$ws.Range("E19").Value = "  -1.95%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.59"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.41%  "

$ws.Range("E21").Value = "  -2.48%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.99"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.29%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.43"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.91%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "233.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.49%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.97"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.10%  "

$ws.Range("E27").Value = "  -4.34%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.89"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.92%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.30"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.05%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "168.77"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.22%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.03"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.67%  "

$ws.Range("E32").Value = "  -2.44%  "

$ws.Range("E34").Value = "  -4.59%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.50"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.65%  "

$ws.Range("E36").Value = "  -5.01%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "16.47"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.72%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0688"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.20%  "

$ws.Range("E39").Value = "  -4.35%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0987"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.19%  "

$ws.Range("E41").Value = "  -3.60%  "

$ws.Range("E42").Value = "  -7.27%  "

$ws.Range("E43").Value = "  -6.06%  "

$ws.Range("D44").Value = "1.959.22"
$ws.Range("E44").Value = "  -3.17%  "

$ws.Range("E45").Value = "  -2.61%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.38"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -8.40%  "

$ws.Range("E47").Value = "  -6.51%  "

$ws.Range("E48").Value = "  -5.35%  "

$ws.Range("D49").Value = "2.494.23"
$ws.Range("E49").Value = "  -2.44%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "51.89"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.40%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.53"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.87%  "
